# Apply updated "dSF" (column F) values for 2021 parker_blake save data.
# These values represent a recalculated/repulled mean-based metric for
# specific rows (column F = "dSF"), per commit message:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$updates = @{
    2  = 3
    3  = 0
    4  = -4
    12 = -3
    14 = -3
    16 = -2
    24 = -17
    25 = 6
    27 = -4
    29 = -1
    32 = -6
    33 = -2
    38 = -3
    44 = 4
    45 = -4
    50 = -4
    51 = -6
    55 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
